$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44769
$ws.Range("J2").Value2 = 50
$ws.Range("K2").Value2 = 20000
$ws.Range("L2").Value2 = 20000
$ws.Range("M2").Value2 = 20000
$ws.Range("P2").Value2 = 1333

$ws.Range("D3").Value2 = 44845
$ws.Range("J3").Value2 = 20
$ws.Range("K3").Value2 = 16000
$ws.Range("L3").Value2 = 16000
$ws.Range("M3").Value2 = 16000
$ws.Range("P3").Value2 = 1067

$ws.Range("D4").Value2 = 44826
$ws.Range("J4").Value2 = 50
$ws.Range("K4").Value2 = 20000
$ws.Range("L4").Value2 = 20000
$ws.Range("M4").Value2 = 20000
$ws.Range("P4").Value2 = 1333

$ws.Range("D5").Value2 = 44819
$ws.Range("J5").Value2 = 100
$ws.Range("K5").Value2 = 20000
$ws.Range("L5").Value2 = 20000
$ws.Range("M5").Value2 = 20000
$ws.Range("P5").Value2 = 1333

$ws.Range("D6").Value2 = 44813
$ws.Range("J6").Value2 = 20
$ws.Range("K6").Value2 = 20000
$ws.Range("L6").Value2 = 20000
$ws.Range("M6").Value2 = 20000
$ws.Range("P6").Value2 = 1333

$ws.Range("D7").Value2 = 44841
$ws.Range("J7").Value2 = 20
$ws.Range("K7").Value2 = 16000
$ws.Range("L7").Value2 = 16000
$ws.Range("M7").Value2 = 16000
$ws.Range("P7").Value2 = 1067

$ws.Range("D8").Value2 = 44830
$ws.Range("J8").Value2 = 25
$ws.Range("K8").Value2 = 12000
$ws.Range("L8").Value2 = 12000
$ws.Range("M8").Value2 = 12000
$ws.Range("P8").Value2 = 800

$ws.Range("D9").Value2 = 44838
$ws.Range("J9").Value2 = 10
$ws.Range("K9").Value2 = 20000
$ws.Range("L9").Value2 = 20000
$ws.Range("M9").Value2 = 20000
$ws.Range("P9").Value2 = 1333

$ws.Range("D10").Value2 = 44755
$ws.Range("J10").Value2 = 50
$ws.Range("K10").Value2 = 20000
$ws.Range("L10").Value2 = 20000
$ws.Range("M10").Value2 = 20000
$ws.Range("P10").Value2 = 1333

$ws.Range("D11").Value2 = 44508
$ws.Range("J11").Value2 = 40
$ws.Range("K11").Value2 = 10000
$ws.Range("L11").Value2 = 10000
$ws.Range("M11").Value2 = 10000
$ws.Range("P11").Value2 = 667

$ws.Range("D12").Value2 = 44525
$ws.Range("J12").Value2 = 40
$ws.Range("K12").Value2 = 8000
$ws.Range("L12").Value2 = 8000
$ws.Range("M12").Value2 = 8000
$ws.Range("P12").Value2 = 533

$ws.Range("D13").Value2 = 44827
$ws.Range("J13").Value2 = 20
$ws.Range("K13").Value2 = 20000
$ws.Range("L13").Value2 = 20000
$ws.Range("M13").Value2 = 20000
$ws.Range("P13").Value2 = 1333

$ws.Range("D14").Value2 = 44771
$ws.Range("J14").Value2 = 40
$ws.Range("K14").Value2 = 20000
$ws.Range("L14").Value2 = 20000
$ws.Range("M14").Value2 = 20000
$ws.Range("P14").Value2 = 1333

$ws.Range("D15").Value2 = 44518
$ws.Range("J15").Value2 = 50
$ws.Range("K15").Value2 = 10000
$ws.Range("L15").Value2 = 10000
$ws.Range("M15").Value2 = 10000
$ws.Range("P15").Value2 = 667

$ws.Range("D16").Value2 = 44839
$ws.Range("J16").Value2 = 80
$ws.Range("K16").Value2 = 16000
$ws.Range("L16").Value2 = 16000
$ws.Range("M16").Value2 = 16000
$ws.Range("P16").Value2 = 1067

$ws.Range("D17").Value2 = 44749
$ws.Range("J17").Value2 = 50
$ws.Range("K17").Value2 = 20000
$ws.Range("L17").Value2 = 20000
$ws.Range("M17").Value2 = 20000
$ws.Range("P17").Value2 = 1333

$ws.Range("D18").Value2 = 45134
$ws.Range("J18").Value2 = 5
$ws.Range("K18").Value2 = 20000
$ws.Range("L18").Value2 = 20000
$ws.Range("M18").Value2 = 20000
$ws.Range("P18").Value2 = 1333

$ws.Range("D19").Value2 = 44812
$ws.Range("J19").Value2 = 80
$ws.Range("K19").Value2 = 20000
$ws.Range("L19").Value2 = 20000
$ws.Range("M19").Value2 = 20000
$ws.Range("P19").Value2 = 1333

$ws.Range("D20").Value2 = 44756
$ws.Range("J20").Value2 = 80
$ws.Range("K20").Value2 = 20000
$ws.Range("L20").Value2 = 20000
$ws.Range("M20").Value2 = 20000
$ws.Range("P20").Value2 = 1333

$ws.Range("D21").Value2 = 45225
$ws.Range("J21").Value2 = 80
$ws.Range("K21").Value2 = 20000
$ws.Range("L21").Value2 = 20000
$ws.Range("M21").Value2 = 20000
$ws.Range("P21").Value2 = 1333

$ws.Range("D22").Value2 = 44837
$ws.Range("J22").Value2 = 80
$ws.Range("K22").Value2 = 16000
$ws.Range("L22").Value2 = 16000
$ws.Range("M22").Value2 = 16000
$ws.Range("P22").Value2 = 1067

$ws.Range("D23").Value2 = 44811
$ws.Range("J23").Value2 = 30
$ws.Range("K23").Value2 = 20000
$ws.Range("L23").Value2 = 20000
$ws.Range("M23").Value2 = 20000
$ws.Range("P23").Value2 = 1333

$ws.Range("D24").Value2 = 44757
$ws.Range("J24").Value2 = 30
$ws.Range("K24").Value2 = 20000
$ws.Range("L24").Value2 = 20000
$ws.Range("M24").Value2 = 20000
$ws.Range("P24").Value2 = 1333

$ws.Range("D25").Value2 = 44767
$ws.Range("J25").Value2 = 50
$ws.Range("K25").Value2 = 20000
$ws.Range("L25").Value2 = 20000
$ws.Range("M25").Value2 = 20000
$ws.Range("P25").Value2 = 1333

$ws.Range("D26").Value2 = 44825
$ws.Range("J26").Value2 = 30
$ws.Range("K26").Value2 = 20000
$ws.Range("L26").Value2 = 20000
$ws.Range("M26").Value2 = 20000
$ws.Range("P26").Value2 = 1333

$ws.Range("D27").Value2 = 44776
$ws.Range("J27").Value2 = 80
$ws.Range("K27").Value2 = 20000
$ws.Range("L27").Value2 = 20000
$ws.Range("M27").Value2 = 20000
$ws.Range("P27").Value2 = 1333

$ws.Range("D28").Value2 = 44824
$ws.Range("J28").Value2 = 20
$ws.Range("K28").Value2 = 20000
$ws.Range("L28").Value2 = 20000
$ws.Range("M28").Value2 = 20000
$ws.Range("P28").Value2 = 1333
